$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text need an explicit Text format first,
# otherwise Excel auto-converts the typed string into a Number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "37.785.81"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "2.044.99"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "227.62"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6
$ws.Range("E6").Value = "  -0.57%  "

# Row 7
$ws.Range("D7").Value = "60.28"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("E9").Value = "  -2.32%  "

# Row 10
$ws.Range("D10").Value = "0.0841"
$ws.Range("E10").Value = "  +3.38%  "

# Row 11
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("D12").Value = "2.346.22"
$ws.Range("E12").Value = "  +0.59%  "

# Row 13
$ws.Range("D13").Value = "14.35"
$ws.Range("E13").Value = "  -1.59%  "

# Row 14
$ws.Range("D14").Value = "21.27"
$ws.Range("E14").Value = "  +1.51%  "

# Row 15
$ws.Range("D15").Value = "5.46"
$ws.Range("E15").Value = "  +5.49%  "

# Row 16
$ws.Range("D16").Value = "0.764"
$ws.Range("E16").Value = "  +0.92%  "

# Row 17
$ws.Range("D17").Value = "2.045.13"
$ws.Range("E17").Value = "  +1.17%  "

# Row 18
$ws.Range("D18").Value = "37.743.28"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19
$ws.Range("D19").Value = "5.94"
$ws.Range("E19").Value = "  -1.56%  "

# Row 20
$ws.Range("D20").Value = "69.39"
$ws.Range("E20").Value = "  -0.43%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  +0.64%  "

# Row 22
$ws.Range("D22").Value = "223.52"
$ws.Range("E22").Value = "  -0.55%  "

# Row 24
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("E25").Value = "  +3.03%  "

# Row 26
$ws.Range("D26").Value = "169.00"
$ws.Range("E26").Value = "  +2.27%  "

# Row 27
$ws.Range("D27").Value = "9.34"
$ws.Range("E27").Value = "  +1.09%  "

# Row 28
$ws.Range("E28").Value = "  -0.56%  "

# Row 29
$ws.Range("D29").Value = "18.77"
$ws.Range("E29").Value = "  -0.74%  "

# Row 30
$ws.Range("E30").Value = "  -0.12%  "

# Row 31
$ws.Range("E31").Value = "  -0.64%  "

# Row 32
$ws.Range("D32").Value = "2.26"
$ws.Range("E32").Value = "  +8.95%  "

# Row 33
$ws.Range("D33").Value = "4.37"
$ws.Range("E33").Value = "  -1.37%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "4.49"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0602"
$ws.Range("E35").Value = "  +0.24%  "

# Row 36
$ws.Range("D36").Value = "6.53"
$ws.Range("E36").Value = "  +3.67%  "

# Row 37
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  +4.50%  "

# Row 38
$ws.Range("D38").Value = "3.47"
$ws.Range("E38").Value = "  +6.68%  "

# Row 39
$ws.Range("E39").Value = "  -0.10%  "

# Row 40
$ws.Range("D40").Value = "18.00"
$ws.Range("E40").Value = "  +6.95%  "

# Row 41
$ws.Range("D41").Value = "1.536.83"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("D42").Value = "97.90"
$ws.Range("E42").Value = "  +1.23%  "

# Row 43
$ws.Range("E43").Value = "  -0.84%  "

# Row 44
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("E45").Value = "  -1.60%  "

# Row 46
$ws.Range("E46").Value = "  +5.98%  "

# Row 47
$ws.Range("E47").Value = "  -0.26%  "

# Row 48
$ws.Range("E48").Value = "  +0.33%  "

# Row 49
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("D50").Value = "7.02"
$ws.Range("E50").Value = "  -1.72%  "

# Row 51
$ws.Range("D51").Value = "2.233.04"
$ws.Range("E51").Value = "  +0.49%  "
